# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# Sheet "展览" (exhibitions)
$wsExpo = $excel.ActiveWorkbook.Worksheets.Item("展览")
$wsExpo.Range("F7").Value = 322
$wsExpo.Range("F9").Value = 4591
$wsExpo.Range("F16").Value = 1889
$wsExpo.Range("F27").Value = 2469
$wsExpo.Range("F28").Value = 1010
$wsExpo.Range("F29").Value = 2441
$wsExpo.Range("F31").Value = 1324

# Sheet "全部类型" (all types) mirrors the same events on different rows
$wsAll = $excel.ActiveWorkbook.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 322
$wsAll.Range("F10").Value = 4591
$wsAll.Range("F28").Value = 2469
$wsAll.Range("F31").Value = 1010
$wsAll.Range("F33").Value = 2441
$wsAll.Range("F34").Value = 1324
